$d = $word.ActiveDocument

# Locate the "Methodology" heading using Find, then identify the paragraph that
# immediately follows it. That paragraph currently holds just a tab character,
# the _GoBack bookmark, and the trailing page break; it is where the new
# Methodology narrative, the blank spacer line, and the "Results" heading need
# to be inserted (ahead of the existing bookmark/page-break).
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Methodology", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Methodology' heading"
}

$headingStart = $searchRange.Start
$headingEnd = $searchRange.End

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Start -le $headingStart -and $pp.Range.End -ge $headingEnd) {
        $targetIndex = $i + 1
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the paragraph following the 'Methodology' heading"
}

$target = $d.Paragraphs.Item($targetIndex)
$r = $target.Range

# Sanity check: this paragraph should be the short tab/bookmark/page-break
# paragraph that currently sits right under the "Methodology" heading.
if ($r.Text.Trim().Length -gt 1) {
    throw "Unexpected content in the target paragraph; aborting to avoid corrupting the document"
}

$xmlPayload = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">The </w:t></w:r><w:r><w:t xml:space="preserve">program operates in three steps that continually repeat in a loop. First, </w:t></w:r><w:r><w:t>the program scans for new messages in its groups and partitions these messages into regular messages and search queries. Secondly, the program inserts the regular messages into the instance of Elasticsearch (the alternative to Solr selected for the project). Lastly, the program responds to search queries by retrieving relevant messages from Elasticsearch and displaying them in the chat.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:tab/><w:t xml:space="preserve"> The program scans for new messages by first retrieving all of the groups in which the Search Bot is a member. If there are any new groups, the program retrieves all of the messages in the new group chat. For the rest of the groups, the program retrieves all of the messages which were posted after the most recent message</w:t></w:r><w:r><w:t xml:space="preserve"> indexed by the program.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t>The messages retrieved by the program are then partitioned into regular messages and search queries. The criteria for a search query is that it directly mentions the account name of the GroupMe Search Bot. The rest of the messages are indexed into Elasticsearch, with key information including the message group, ID, sender, text, and timestamp being preserved.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:tab/><w:t>The searches are then responded to by extracting the space-separated keywords from the search text and querying the Elasticsearch instance for messages that match these keywords. Note that only messages from the group being searched are considered, even though the Elasticsearch instance stores messages from multiple groups at once. The program then sorts all matches by timestamp, with the most recent messages being displayed first. A limit of ten messages are displayed such that the bot does not spam the chat with unnecessary information. Lastly, the program constructs a message which returns each match on a separate line. If the search result takes up more than the 1000-character message limit, the response is broken apart into multiple messages.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Results</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:tab/></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/></w:pPr><w:r><w:br w:type="page"/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r.InsertXML($xmlPayload)

Write-Host "Inserted Methodology body text and the Results heading."
